$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New lat/lon/sitenumber data (sitename column removed)
$data = @(
    @(30.977402, -83.36899699999999, 1),
    @(32.515813, -86.377325, 2),
    @(42.23498, -88.30540999999999, 3),
    @(33.870013, -118.377777, 4),
    @(34.014929, -118.205387, 5),
    @(40.731099, -74.173067, 6),
    @(37.81144, -121.29348, 7),
    @(44.85387, -93.04713, 8),
    @(41.18661, -111.94904, 9),
    @(40.71239, -74.5847, 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove column D (sitename) entirely, including header
$ws.Range("D1:D11").Delete()
